$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.200234174728394
$ws.Range("B1").Value = 1.877100467681885
$ws.Range("C1").Value = 4.516806125640869
$ws.Range("D1").Value = 1.578817009925842
$ws.Range("E1").Value = 0.5272939801216125
